$wb = $excel.ActiveWorkbook

# Reorder worksheet tabs:
#   Observables and Functions move to just before Reactions.
#   Stop conditions moves to just before References.
$wb.Worksheets.Item("Observables").Move($wb.Worksheets.Item("Reactions"))
$wb.Worksheets.Item("Functions").Move($wb.Worksheets.Item("Reactions"))
$wb.Worksheets.Item("Stop conditions").Move($wb.Worksheets.Item("References"))

# Populate header row for the "Observables" sheet
$wsObservables = $wb.Worksheets.Item("Observables")
$wsObservables.Range("A1").Value = "Id"
$wsObservables.Range("B1").Value = "Name"
$wsObservables.Range("C1").Value = "Model"
$wsObservables.Range("D1").Value = "Species"
$wsObservables.Range("E1").Value = "Observables"
$wsObservables.Range("F1").Value = "Comments"
$wsObservables.Range("A1:F1").Select()

# Populate header row for the "Functions" sheet
$wsFunctions = $wb.Worksheets.Item("Functions")
$wsFunctions.Range("A1").Value = "Id"
$wsFunctions.Range("B1").Value = "Name"
$wsFunctions.Range("C1").Value = "Model"
$wsFunctions.Range("D1").Value = "Expression"
$wsFunctions.Range("E1").Value = "Comments"
$wsFunctions.Range("A1:E1").Select()

# Populate header row for the "Stop conditions" sheet
$wsStopConditions = $wb.Worksheets.Item("Stop conditions")
$wsStopConditions.Range("A1").Value = "Id"
$wsStopConditions.Range("B1").Value = "Name"
$wsStopConditions.Range("C1").Value = "Model"
$wsStopConditions.Range("D1").Value = "Expression"
$wsStopConditions.Range("E1").Value = "Comments"
$wsStopConditions.Range("A1:E1").Select()
$wsStopConditions.Activate()
